$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2023 -2024" / "JURY DU 14 MARS 2024" header labels in F4/F5 are removed
# (template reverted to a date-agnostic model: ModeleS3Jury)
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()

# With the large header text gone, rows 4 & 5 shrink to their natural height
$ws.Rows.Item(4).RowHeight = 19.7
$ws.Rows.Item(5).RowHeight = 19.7

# Touch the bottom-right corner of the (already present, empty) data rows so the
# sheet's extent properly covers every row down to 57, as in the reference template
$ws.Cells.Item(57, 1).NumberFormat = $ws.Cells.Item(57, 1).NumberFormat

# Move/save the active selection on F4 (reflecting where the edit took place)
$ws.Range("F4").Select() | Out-Null
